$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(113).Insert()

$ws.Cells.Item(113, 1).Value = 7
$ws.Cells.Item(113, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(113, 3).Value = "Ñuble"
$ws.Cells.Item(113, 4).Value = 44455
$ws.Cells.Item(113, 4).NumberFormat = $ws.Cells.Item(114, 4).NumberFormat
$ws.Cells.Item(113, 5).Value = 16
$ws.Cells.Item(113, 6).Value = 100114013
$ws.Cells.Item(113, 7).Value = "Zanahoria"
$ws.Cells.Item(113, 8).Value = "Sin especificar"
$ws.Cells.Item(113, 9).Value = "Primera"
$ws.Cells.Item(113, 10).Value = 300
$ws.Cells.Item(113, 11).Value = 7500
$ws.Cells.Item(113, 12).Value = 8000
$ws.Cells.Item(113, 13).Value = 7750
$ws.Cells.Item(113, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(113, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(113, 16).Value = 388
$ws.Cells.Item(113, 17).Value = 20
$ws.Cells.Item(113, 18).Value = "Hortaliza"
